$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new master data rows for "Registration Acknowledgement Template - Part 4"
$ws.Cells.Item(122,1).Value = "reg-ack-template-part4"
$ws.Cells.Item(122,2).Value = "Registration Acknowledgement Template - Part 4"
$ws.Cells.Item(122,3).Value = "eng"
$ws.Cells.Item(122,4).Value = $true
$ws.Cells.Item(122,5).Value = "superadmin"
$ws.Cells.Item(122,6).Value = "now()"

$ws.Cells.Item(123,1).Value = "reg-ack-template-part4"
$ws.Cells.Item(123,2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(123,3).Value = "ara"
$ws.Cells.Item(123,4).Value = $true
$ws.Cells.Item(123,5).Value = "superadmin"
$ws.Cells.Item(123,6).Value = "now()"

$ws.Cells.Item(124,1).Value = "reg-ack-template-part4"
$ws.Cells.Item(124,2).Value = "accusé de réception"
$ws.Cells.Item(124,3).Value = "fra"
$ws.Cells.Item(124,4).Value = $true
$ws.Cells.Item(124,5).Value = "superadmin"
$ws.Cells.Item(124,6).Value = "now()"

# Move the selection below the newly added data, matching the author's saved view
$ws.Range("A125:XFD1048576").Select()
